# Generate Report for Handback
#
# This applies the "handback" refresh to the localization-status workbook:
#   - the "In Translation" status becomes "Handed back: in sync with en-US"
#     everywhere it is used (Overview + per-locale sheets)
#   - per-locale sheets (zh-cn / de-de) get their "Latest Target File" /
#     "Latest Handback File" columns populated with the source .md files
#     (with hyperlinks, matching column A) and handback file names
#   - handback/generate timestamps are refreshed
#   - a handful of columns are widened to fit the new, longer content

$wb = $excel.ActiveWorkbook

$ovWs   = $wb.Worksheets.Item("Overview")
$zhWs   = $wb.Worksheets.Item("zh-cn")
$deWs   = $wb.Worksheets.Item("de-de")

$mdFile1 = "454861c3-9191-4444-9262-2a3b4dda91cc.md"
$mdFile2 = "b1658cf1-a484-4650-b714-88d69cf0de51.md"
$mdUrl1  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/969c38594985b2a435e8baa67e4d4b0dd3bc6ab2/e2e/$mdFile1"
$mdUrl2  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/969c38594985b2a435e8baa67e4d4b0dd3bc6ab2/e2e/$mdFile2"

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Status text: "In Translation" -> "Handed back: in sync with en-US"
#    (every cell that currently shows the old status)
# ---------------------------------------------------------------------
$ovWs.Range("E2").Value = $newStatus
$ovWs.Range("F2").Value = $newStatus
$ovWs.Range("E3").Value = $newStatus
$ovWs.Range("F3").Value = $newStatus

$zhWs.Range("C2").Value = $newStatus
$zhWs.Range("C3").Value = $newStatus

$deWs.Range("C2").Value = $newStatus
$deWs.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: fill in the target/handback file columns + hyperlinks,
#    refresh the handback datetime
# ---------------------------------------------------------------------
$zhWs.Range("I2").Value = $mdFile1
$zhWs.Hyperlinks.Add($zhWs.Range("I2"), $mdUrl1, "", "", $mdFile1) | Out-Null
$zhWs.Range("J2").Value = "454861c3-9191-4444-9262-2a3b4dda91cc.d9fcfb31ff4a1170db7d58ed47c24e6694ac3b5a.zh-cn.xlf"
$zhWs.Range("K2").Value = "2016-08-30 20:26:54"

$zhWs.Range("I3").Value = $mdFile2
$zhWs.Hyperlinks.Add($zhWs.Range("I3"), $mdUrl2, "", "", $mdFile2) | Out-Null
$zhWs.Range("J3").Value = "b1658cf1-a484-4650-b714-88d69cf0de51.ceaaf243d8f2b640ab9bd342bf340065bc76bdb1.zh-cn.xlf"
$zhWs.Range("K3").Value = "2016-08-30 20:26:54"

# ---------------------------------------------------------------------
# 3. de-de sheet: same treatment, plus the handoff-file / generate-date
#    columns refresh
# ---------------------------------------------------------------------
$deWs.Range("G2").Value = "2016-08-30 20:27:04"
$deWs.Range("I2").Value = $mdFile1
$deWs.Hyperlinks.Add($deWs.Range("I2"), $mdUrl1, "", "", $mdFile1) | Out-Null
$deWs.Range("J2").Value = "454861c3-9191-4444-9262-2a3b4dda91cc.d9fcfb31ff4a1170db7d58ed47c24e6694ac3b5a.de-de.xlf"
$deWs.Range("K2").Value = "b1658cf1-a484-4650-b714-88d69cf0de51.ceaaf243d8f2b640ab9bd342bf340065bc76bdb1.de-de.xlf"

$deWs.Range("G3").Value = "2016-08-30 20:27:04"
$deWs.Range("I3").Value = $mdFile2
$deWs.Hyperlinks.Add($deWs.Range("I3"), $mdUrl2, "", "", $mdFile2) | Out-Null
$deWs.Range("J3").Value = "2016-08-30 20:27:04"
$deWs.Range("K3").Value = "b1658cf1-a484-4650-b714-88d69cf0de51.ceaaf243d8f2b640ab9bd342bf340065bc76bdb1.de-de.xlf"

# ---------------------------------------------------------------------
# 4. Widen columns to fit the newly-populated / longer content
#    (ColumnWidth is in characters; stored sheet width = ColumnWidth + 5/6)
# ---------------------------------------------------------------------
$wide30 = 29.144371396019366   # -> stored width ~29.98 (longest status text)
$wide40 = 39.166666666666664   # -> stored width 40 (matches other "40" columns)

$ovWs.Columns.Item(5).ColumnWidth = $wide30   # Overview!E (zh-cn)
$ovWs.Columns.Item(6).ColumnWidth = $wide30   # Overview!F (de-de)

$zhWs.Columns.Item(3).ColumnWidth  = $wide30  # zh-cn!C  (Status)
$zhWs.Columns.Item(9).ColumnWidth  = $wide40  # zh-cn!I  (Latest Target File)
$zhWs.Columns.Item(10).ColumnWidth = $wide40  # zh-cn!J  (Latest Handback File)

$deWs.Columns.Item(3).ColumnWidth  = $wide30  # de-de!C  (Status)
$deWs.Columns.Item(9).ColumnWidth  = $wide40  # de-de!I  (Latest Target File)
$deWs.Columns.Item(10).ColumnWidth = $wide40  # de-de!J  (Latest Handback File)
